# Apply updated per-subject statistics after branching/editing the program
# to disregard certain variables.

$wb = $excel.ActiveWorkbook

# --- "Program Control" sheet: Subroutine count picks up for subjects 21 & 23 ---
$wsControl = $wb.Worksheets.Item("Program Control")
$wsControl.Range("B22").Value = 1
$wsControl.Range("B24").Value = 1

# --- "Variables" sheet: updated variable-tracking stats after the edit ---
$wsVars = $wb.Worksheets.Item("Variables")

$wsVars.Range("B7").Value = 12
$wsVars.Range("E7").Value = 5
$wsVars.Range("G7").Value = 8.916666666666666

$wsVars.Range("B8").Value = 13
$wsVars.Range("F8").Value = 2
$wsVars.Range("G8").Value = 9

$wsVars.Range("H9").Value = 0

$wsVars.Range("B12").Value = 10
$wsVars.Range("F12").Value = 2
$wsVars.Range("G12").Value = 4.8

$wsVars.Range("F15").Value = 2

$wsVars.Range("B16").Value = 16
$wsVars.Range("E16").Value = 5
$wsVars.Range("G16").Value = 6

$wsVars.Range("B19").Value = 18
$wsVars.Range("G19").Value = 6.277777777777778

$wsVars.Range("B21").Value = 16
$wsVars.Range("C21").Value = 6
$wsVars.Range("E21").Value = 7
$wsVars.Range("G21").Value = 5.125

$wsVars.Range("B24").Value = 17
$wsVars.Range("E24").Value = 8
$wsVars.Range("G24").Value = 5.352941176470588

$wsVars.Range("B28").Value = 23
$wsVars.Range("E28").Value = 7
$wsVars.Range("G28").Value = 5.434782608695652
